# Append four new "Positioning" target rows (rows 30-33) to the
# external_targets sheet, mirroring the existing row layout/typing, and
# move the sheet selection to D10 (matching the updated view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-TargetRow($Row, $Pk, $OperationTypeFk, $Level2Fk, $KeyJson, $DataJson) {
    $ws.Cells.Item($Row, 1).Value = $Pk
    $ws.Cells.Item($Row, 2).Value = $OperationTypeFk
    $ws.Cells.Item($Row, 3).Value = $Level2Fk
    $ws.Cells.Item($Row, 4).Value = $KeyJson
    $ws.Cells.Item($Row, 5).Value = $DataJson

    # start_date / received_time columns are stored as text (matching the
    # "@" number format already used by the existing rows in this sheet).
    $ws.Cells.Item($Row, 6).NumberFormat = "@"
    $ws.Cells.Item($Row, 6).Value = "2019-04-01 00:00:00"

    $ws.Cells.Item($Row, 7).Value = "NULL"

    $ws.Cells.Item($Row, 8).NumberFormat = "@"
    $ws.Cells.Item($Row, 8).Value = "2019-04-03 09:38:27"

    $ws.Cells.Item($Row, 9).Value = "Positioning"
}

Add-TargetRow 30 58 3 319 '{"Group Name": "Pringles_FTT_Tubes"}' '{"Target": 90, "Value 1": 189, "Value 2": "Fun times together Tubes", "Value 3": "", "Parameter 1": "brand_fk", "Parameter 2": "PDH Format", "Parameter 3": ""}'

Add-TargetRow 31 59 3 319 '{"Group Name": "Hula Hoops_LMP_Snacks"}' '{"Target": 90, "Value 1": 138, "Value 2": "LMP Snacks", "Value 3": "", "Parameter 1": "brand_fk", "Parameter 2": "PDH Format", "Parameter 3": ""}'

Add-TargetRow 32 60 3 316 '{"Group Name": "DORITOS GROUP"}' '{"Target": 0, "Value 1": 136, "Value 2": "", "Value 3": "", "Parameter 1": "brand_fk", "Parameter 2": "", "Parameter 3": ""}'

Add-TargetRow 33 61 3 316 '{"Group Name": "Walkers Crisps_Small MP PC"}' '{"Target": 0, "Value 1": 199, "Value 2": 7, "Value 3": "SMP PC", "Parameter 1": "brand_fk", "Parameter 2": "sub_category_fk", "Parameter 3": "PDH Sub-segment"}'

# Move the active selection to D10 and scroll the view back to the top-left
# (A1), matching the refreshed window state recorded for the sheet.
$ws.Range("A1").Select() | Out-Null
$ws.Range("D10").Select() | Out-Null

# Minor cosmetic tab-ratio tweak recorded alongside the data edit.
$excel.ActiveWindow.TabRatio = 0.993
